$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meter Entries")

$rows = 7..18
foreach ($r in $rows) {
    $c = $ws.Range("L$r")
    $c.NumberFormat = "@"
    $c.HorizontalAlignment = 1
    $c.Value = "Not Available"
}
